$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.624.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.296.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.289.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("E9").Value = "  -3.33%  "
$ws.Range("E10").Value = "  -6.69%  "
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.16%  "
$ws.Range("E13").Value = "  -3.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.822.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "570.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.551.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.291.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.27%  "
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.14%  "
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "555.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.747.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -3.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "32.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -4.41%  "
$ws.Range("E42").Value = "  -9.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0666"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.02%  "
$ws.Range("E45").Value = "  -6.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.327"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -12.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  -3.42%  "
$ws.Range("E51").Value = "  -4.81%  "
